$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing bold/centered/bordered format from A16 onto the new A17:A19 cells
# (and B17:B19 to match the plain style used by the existing B column) before writing values,
# so the new rows visually match the rest of the table.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# Row 10: Gaussian-Quadrature
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.004001539351017
$ws.Range("D10").Value = 0.9977992842888757
$ws.Range("E10").Value = 0.9969934576433873
$ws.Range("F10").Value = 0.9965505781855611
$ws.Range("G10").Value = 1.004001539351017
$ws.Range("H10").Value = 0.9977992842888757
$ws.Range("I10").Value = 1.002890434648806
$ws.Range("J10").Value = 1.004183012944848
$ws.Range("K10").Value = 1.001764705882353
$ws.Range("L10").Value = 0.9988381534760072
$ws.Range("M10").Value = 1.004001539351017
$ws.Range("N10").Value = 0.9973963709661315
$ws.Range("O10").Value = 0.9988362148672103
$ws.Range("P10").Value = 1.000377645802607

# Row 11: Spiral-90deg-10rot-5space
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9883836777873508
$ws.Range("D11").Value = 1.022357172799638
$ws.Range("E11").Value = 0.9969374541203054
$ws.Range("F11").Value = 1.006233435133741
$ws.Range("G11").Value = 0.9883836777873508
$ws.Range("H11").Value = 1.022357172799638
$ws.Range("I11").Value = 0.9924249541335471
$ws.Range("J11").Value = 1.005203906099783
$ws.Range("K11").Value = 0.9959337207146849
$ws.Range("L11").Value = 1.015525468281536
$ws.Range("M11").Value = 0.9883836777873508
$ws.Range("N11").Value = 1.009647313459972
$ws.Range("O11").Value = 1.003477934960259
$ws.Range("P11").Value = 1.002874973633823

# Row 12: Spiral-90deg-15rot-5space
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9883562099745855
$ws.Range("D12").Value = 1.022470091483814
$ws.Range("E12").Value = 0.9969643435410405
$ws.Range("F12").Value = 1.006281656054412
$ws.Range("G12").Value = 0.9883562099745855
$ws.Range("H12").Value = 1.022470091483814
$ws.Range("I12").Value = 0.9923321016573046
$ws.Range("J12").Value = 1.005139870111546
$ws.Range("K12").Value = 0.9959242338054352
$ws.Range("L12").Value = 1.01559592239892
$ws.Range("M12").Value = 0.9883562099745855
$ws.Range("N12").Value = 1.009717217512427
$ws.Range("O12").Value = 1.003518075263463
$ws.Range("P12").Value = 1.002883053628382

# Row 13: Spiral-90deg-10rot-3space
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9883812341737249
$ws.Range("D13").Value = 1.022381054356469
$ws.Range("E13").Value = 0.9969516242356437
$ws.Range("F13").Value = 1.006248178022062
$ws.Range("G13").Value = 0.9883812341737249
$ws.Range("H13").Value = 1.022381054356469
$ws.Range("I13").Value = 0.9923862917060815
$ws.Range("J13").Value = 1.005165588612935
$ws.Range("K13").Value = 0.995933012218772
$ws.Range("L13").Value = 1.015542332008813
$ws.Range("M13").Value = 0.9883812341737249
$ws.Range("N13").Value = 1.009666339296056
$ws.Range("O13").Value = 1.003490522696975
$ws.Range("P13").Value = 1.002873664416813

# Row 14: NoRotation-tilt60deg
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.9925039999999995
$ws.Range("D14").Value = 1.003736000000001
$ws.Range("E14").Value = 0.9920800000000005
$ws.Range("F14").Value = 0.9983280000000007
$ws.Range("G14").Value = 0.9925039999999995
$ws.Range("H14").Value = 1.003736000000001
$ws.Range("I14").Value = 1.005219999999999
$ws.Range("J14").Value = 1.013131999999999
$ws.Range("K14").Value = 0.9972920000000002
$ws.Range("L14").Value = 1.004780000000002
$ws.Range("M14").Value = 0.9925039999999995
$ws.Range("N14").Value = 0.9979080000000007
$ws.Range("O14").Value = 0.9966620000000004
$ws.Range("P14").Value = 1.000884

# Row 15: Rotation-NoTilt
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0.98
$ws.Range("E15").Value = 0.99
$ws.Range("F15").Value = 0.99
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0.98
$ws.Range("I15").Value = 1.02
$ws.Range("J15").Value = 1.02
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.99
$ws.Range("M15").Value = 1
$ws.Range("N15").Value = 0.985
$ws.Range("O15").Value = 0.99
$ws.Range("P15").Value = 0.99875

# Row 16: Rotation-60detTilt
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 1.000941521510403
$ws.Range("D16").Value = 0.9888565141503992
$ws.Range("E16").Value = 0.9955855755263997
$ws.Range("F16").Value = 0.9949921193983998
$ws.Range("G16").Value = 1.000941521510403
$ws.Range("H16").Value = 0.9888565141503992
$ws.Range("I16").Value = 1.0112700008448
$ws.Range("J16").Value = 1.010862530559998
$ws.Range("K16").Value = 1.000455529062398
$ws.Range("L16").Value = 0.9948403554303968
$ws.Range("M16").Value = 1.000941521510403
$ws.Range("N16").Value = 0.9922210448383995
$ws.Range("O16").Value = 0.9950939326464003
$ws.Range("P16").Value = 0.9997255183103992

# Row 17: HexGrid-90degTilt5degRes
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 1.00155321537034
$ws.Range("D17").Value = 1.001014174664779
$ws.Range("E17").Value = 1.001585625256139
$ws.Range("F17").Value = 1.000962094148966
$ws.Range("G17").Value = 1.00155321537034
$ws.Range("H17").Value = 1.001014174664779
$ws.Range("I17").Value = 1.000430032913838
$ws.Range("J17").Value = 1.0008019392438
$ws.Range("K17").Value = 1.000806155356935
$ws.Range("L17").Value = 1.00129729386553
$ws.Range("M17").Value = 1.00155321537034
$ws.Range("N17").Value = 1.001299899960459
$ws.Range("O17").Value = 1.001278777360056
$ws.Range("P17").Value = 1.001056316352541

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.000434122083101
$ws.Range("D18").Value = 1.002256999084803
$ws.Range("E18").Value = 1.001620367298374
$ws.Range("F18").Value = 1.000573938224906
$ws.Range("G18").Value = 1.000434122083101
$ws.Range("H18").Value = 1.002256999084803
$ws.Range("I18").Value = 1.00014262908846
$ws.Range("J18").Value = 1.001183795275812
$ws.Range("K18").Value = 1.000902179082475
$ws.Range("L18").Value = 1.002419018510564
$ws.Range("M18").Value = 1.000434122083101
$ws.Range("N18").Value = 1.001938683191589
$ws.Range("O18").Value = 1.001221356672796
$ws.Range("P18").Value = 1.001191631081062

# Row 19: HexGrid-60degTilt5degRes
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 1.000854996607823
$ws.Range("D19").Value = 1.002670177542221
$ws.Range("E19").Value = 1.001170958968433
$ws.Range("F19").Value = 1.001482482453176
$ws.Range("G19").Value = 1.000854996607823
$ws.Range("H19").Value = 1.002670177542221
$ws.Range("I19").Value = 0.9998245457093347
$ws.Range("J19").Value = 1.000811100055283
$ws.Range("K19").Value = 1.000546094034378
$ws.Range("L19").Value = 1.002426579449365
$ws.Range("M19").Value = 1.000854996607823
$ws.Range("N19").Value = 1.001920568255327
$ws.Range("O19").Value = 1.001278777360056
$ws.Range("P19").Value = 1.001223366852502
